$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "250.18"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.83"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.433"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.381"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8151"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9214"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1440"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07479"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03113"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03101"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09356"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.769"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001585"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04780"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005796"

$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006416"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005035"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001032"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001502"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.703"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.181"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3304"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1322"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003003"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04023"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006777"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002713"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008042"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005809"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
